# Implement ACF in more risk groups
# Add urban and rural poor to risk groups that we can screen with ACF.
#
# On the "constants" sheet, insert 10 new parameter rows (urbanpoor + ruralpoor
# xpertacf economics blocks) right above the existing
# "econ_unitcost_engage_lowquality" block, pushing that block (and the final
# saturation row) down by 10 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")
$ws.Activate()

# The "econ_unitcost_engage_lowquality" block currently starts at row 141.
# Insert 10 blank rows above it (formats/styles copied down from row 140,
# matching how the surrounding parameter rows are styled).
$ws.Range("A141:A150").EntireRow.Insert()

$newRows = @(
  @("econ_unitcost_xpertacf_urbanpoor", 30.26),
  @("econ_inflectioncost_xpertacf_urbanpoor", 0),
  @("econ_startupcost_xpertacf_urbanpoor", 662),
  @("econ_startupduration_xpertacf_urbanpoor", 1),
  @("econ_saturation_xpertacf_urbanpoor", 0.9),
  @("econ_unitcost_xpertacf_ruralpoor", 30.26),
  @("econ_inflectioncost_xpertacf_ruralpoor", 0),
  @("econ_startupcost_xpertacf_ruralpoor", 662),
  @("econ_startupduration_xpertacf_ruralpoor", 1),
  @("econ_saturation_xpertacf_ruralpoor", 0.9)
)

$r = 141
foreach ($row in $newRows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $r = $r + 1
}

# Reflect the scroll position / active selection the author ended up with.
try {
  $excel.ActiveWindow.ScrollRow = 40
  $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("E53").Select() | Out-Null
